# Apply the edits described by the commit diff.
#
# The diff is dominated by cosmetic re-splitting of runs (Word inserting
# w:proofErr spellcheck/grammar markers around re-typed text) and by
# watermark shape/anchor id churn that is regenerated every time the file
# is resaved by Word -- none of that is semantically meaningful and none
# of it is reproducible from a generic script, so we only reproduce the
# genuine textual edits:
#
#   1. {{SEXO_5}} -> {{SEXO_9}} in the opening "SER PERSONA FISICA..."
#      declaration of the first party.
#   2. "SOBRE LA ESCRITURACION POR CUESTIONES AJENAS A ESTE" ->
#      "SOBRE LA COMPRAVENTA POR CUESTIONES AJENAS A ESTE" (both
#      occurrences).

$d = $word.ActiveDocument

# 1) {{SEXO_5}} -> {{SEXO_9}} -- only in the first party's declaration
#    (other paragraphs also contain a literal "{{SEXO_5}}" token but with
#    different surrounding text, so anchor the search on enough context
#    to hit only the intended run).
$d.Content.Find.Execute(
    "MEXICAN{{SEXO_5}}, MAYOR DE EDAD, CON DOMICILIO PARA LOS EFECTOS",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "MEXICAN{{SEXO_9}}, MAYOR DE EDAD, CON DOMICILIO PARA LOS EFECTOS",
    2
)

# 2) ESCRITURACIÓN -> COMPRAVENTA (both occurrences, whole document)
$d.Content.Find.Execute(
    "SOBRE LA ESCRITURACIÓN POR CUESTIONES AJENAS A ÉSTE",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SOBRE LA COMPRAVENTA POR CUESTIONES AJENAS A ÉSTE",
    2
)
